$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------
# 1) Rows whose entire content (columns B..AB) must be swapped
#    between two row numbers. Column A (id) is left untouched.
# ---------------------------------------------------------------
$swapRows = @(
  @{Row=85; B=6992623; C='Thailand Premier League'; D=45261.375; E='Ratchaburi FC'; F='Chiangrai Utd'; G=3; H=0; I='H'; J=1.7; K=3.75; L=4.2; M=1.7; N=3.75; O=4.333; P=-0.75; Q=1.925; R=1.875; S=2.5; T=1.85; U=1.95; V=0.7; W=-1; X=-1; Y=0.925; Z=-1; AA=0.8500000000000001; AB=-1}
  @{Row=86; B=6992620; C='Thailand Premier League'; D=45261.375; E='Uthai Thani FC'; F='Sukhothai FC'; G=0; H=0; I='D'; J=1.95; K=3.5; L=3.4; M=2.1; N=3.4; O=3; P=-0.25; Q=1.875; R=1.925; S=2.75; T=1.8; U=2; V=-1; W=2.4; X=-1; Y=-0.5; Z=0.4625; AA=-1; AB=1}
  @{Row=179; B=8026714; C='Thailand Premier League'; D=45385.375; E='BG Pathum United'; F='Buriram United'; G=1; H=1; I='D'; J=3; K=3.6; L=2; M=3.1; N=3.75; O=1.95; P=0.5; Q=1.825; R=1.975; S=2.75; T=1.85; U=1.95; V=-1; W=2.75; X=-1; Y=0.825; Z=-1; AA=-1; AB=0.95}
  @{Row=180; B=6992695; C='Thailand Premier League'; D=45385.375; E='Muang Thong United'; F='Uthai Thani FC'; G=5; H=2; I='H'; J=2.1; K=3.75; L=2.7; M=1.95; N=3.8; O=2.9; P=-0.25; Q=1.8; R=2; S=3; T=1.825; U=1.975; V=0.95; W=-1; X=-1; Y=0.8; Z=-1; AA=0.825; AB=-1}
)

foreach ($r in $swapRows) {
  foreach ($key in $r.Keys) {
    if ($key -eq 'Row') { continue }
    $cellRef = "$key$($r.Row)"
    $ws.Range($cellRef).Value = $r[$key]
  }
}

# ---------------------------------------------------------------
# 2) Brand new rows appended at the bottom of the table (216..223)
# ---------------------------------------------------------------
$newRows = @(
  @{Row=216; A=214; B=6992735; C='Thailand Premier League'; D=45422.375; E='BG Pathum United'; F='Khonkaen United'; G=3; H=2; I='H'; J=1.333; K=5; L=6; M=1.25; N=5.75; O=7; P=-1.75; Q=1.825; R=1.975; S=3.75; T=1.9; U=1.9; V=0.25; W=-1; X=-1; Y=-1; Z=0.9750000000000001; AA=0.8999999999999999; AB=-1}
  @{Row=217; A=215; B=6992732; C='Thailand Premier League'; D=45423.33333333334; E='Ratchaburi FC'; F='Muang Thong United'; G=1; H=2; I='A'; J=2.1; K=3.5; L=2.8; M=2.3; N=3.4; O=2.55; P=0; Q=1.775; R=2.025; S=2.75; T=1.85; U=1.95; V=-1; W=-1; X=1.55; Y=-1; Z=1.025; AA=0.425; AB=-0.5}
  @{Row=218; A=216; B=6992736; C='Thailand Premier League'; D=45423.35416666666; E='Chiangrai Utd'; F='Nakhon Pathom FC'; G=0; H=0; I='D'; J=1.909; K=2.875; L=4.2; M=1.8; N=2.9; O=4.5; P=-0.5; Q=1.875; R=1.925; S=2.5; T=1.975; U=1.825; V=-1; W=1.9; X=-1; Y=-1; Z=0.925; AA=-1; AB=0.825}
  @{Row=219; A=217; B=6992340; C='Thailand Premier League'; D=45423.41666666666; E='Sukhothai FC'; F='Police Tero FC'; G=2; H=0; I='H'; J=1.6; K=4; L=4.2; M=1.533; N=4.2; O=4.5; P=-1; Q=1.85; R=1.95; S=3; T=1.8; U=2; V=0.5329999999999999; W=-1; X=-1; Y=0.8500000000000001; Z=-1; AA=-1; AB=1}
  @{Row=220; A=218; B=6992731; C='Thailand Premier League'; D=45424.3125; E='Chonburi'; F='Prachuap FC'; G=1; H=1; I='D'; J=2.15; K=3.3; L=2.875; M=2.375; N=3.25; O=2.6; P=0; Q=1.775; R=2.025; S=2.75; T=1.925; U=1.875; V=-1; W=2.25; X=-1; Y=0; Z=0; AA=-1; AB=0.875}
  @{Row=221; A=219; B=6992730; C='Thailand Premier League'; D=45424.33333333334; E='Buriram United'; F='Port FC'; G=1; H=1; I='D'; J=1.5; K=4; L=5; M=1.571; N=3.8; O=4.5; P=-1; Q=1.95; R=1.85; S=3; T=1.75; U=1.95; V=-1; W=2.8; X=-1; Y=-1; Z=0.8500000000000001; AA=-1; AB=0.95}
  @{Row=222; A=220; B=6992734; C='Thailand Premier League'; D=45424.375; E='Bangkok United'; F='Trat FC'; G=5; H=0; I='H'; J=1.2; K=5.5; L=9.5; M=1.181; N=6; O=10; P=-2.25; Q=1.95; R=1.85; S=3.75; T=1.95; U=1.85; V=0.181; W=-1; X=-1; Y=0.95; Z=-1; AA=0.95; AB=-1}
  @{Row=223; A=221; B=6992733; C='Thailand Premier League'; D=45425.35416666666; E='Lamphun Warrior FC'; F='Uthai Thani FC'; G=1; H=2; I='A'; J=1.909; K=3; L=3.8; M=1.833; N=3.2; O=4; P=-0.5; Q=1.9; R=1.9; S=3; T=1.875; U=1.925; V=-1; W=-1; X=3; Y=-1; Z=0.8999999999999999; AA=0; AB=0}
)

# Copy the style (formats only) of the last existing data row (215)
# down onto each new row so that column A keeps the bold/border style
# and column D keeps its date number format, matching the rest of
# the table.
$lastRow = 215
foreach ($r in $newRows) {
  $ws.Range("A$lastRow`:AB$lastRow").Copy() | Out-Null
  $ws.Range("A$($r.Row)`:AB$($r.Row)").PasteSpecial(-4122) | Out-Null
}
$excel.CutCopyMode = 0

foreach ($r in $newRows) {
  foreach ($key in $r.Keys) {
    if ($key -eq 'Row') { continue }
    $cellRef = "$key$($r.Row)"
    $ws.Range($cellRef).Value = $r[$key]
  }
}

Write-Host "Edit complete"
